$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers (subject numbers)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) updated meanEMG values
$ws.Range("B2").Value = 281.90691521578503
$ws.Range("C2").Value = 262.85707341393442
$ws.Range("D2").Value = 281.67242896552921
$ws.Range("E2").Value = 260.71443222478274

# Row 3 (STR) updated meanEMG values
$ws.Range("B3").Value = 308.11008337506928
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 313.81384839709028
$ws.Range("E3").Value = 256.34900221840485

# Update selection to reflect the new data extent used
$ws.Range("B1:E3").Select() | Out-Null
